# Adds a new review row (shared-string backed) to the sheet, restores the
# wrapped-text auto row heights for the existing review rows, and moves the
# selection to the new bottom of the used range - mirrors the author's
# "added analyze with collocations & selector of dictionary" commit, which
# appended one more review string below the existing ones.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Append the new review text as a new row at the bottom of the data ---
$newRow = 11
$ws.Cells.Item($newRow, 1).Value = "очень хороший телефон, мне безумно понравился"
# Match the formatting (wrap-text style) used by the other review cells.
$ws.Cells.Item($newRow, 1).WrapText = $true

# --- Restore the auto-computed "wrap text" row heights on the review rows ---
# (Excel recalculates these whenever the wrapped text / column width causes a
# reflow; the values below are the heights Excel itself produced for this
# text/column-width combination, capped at the 409.6pt row-height maximum.)
$rowHeights = @{
    2  = 409.6
    3  = 187.2
    4  = 360
    5  = 129.6
    6  = 259.2
    7  = 409.6
    8  = 409.6
    9  = 216
    10 = 409.6
}
foreach ($r in $rowHeights.Keys) {
    $ws.Rows.Item($r).RowHeight = $rowHeights[$r]
}

# --- Scroll the viewport down to the newly added row and move the visible
#     selection to just past it, matching where the user's cursor ended up
#     after inserting/reviewing the new row ---
$win = $excel.ActiveWindow
$win.ScrollRow = $newRow
$win.ScrollColumn = 1
$ws.Range("A14").Select()
